# New weekly price observation for Piña (Macroferia Regional de Talca).
# Insert a new data row at row 250, pushing the existing rows 250-264 down
# to 251-265 (dimension grows from A1:T264 to A1:T265), then populate the
# new row with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(250).Insert()

$ws.Cells.Item(250, 1).Value = 5
$ws.Cells.Item(250, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(250, 3).Value = 'Maule'
$ws.Cells.Item(250, 4).Value = '2022-07-04'
$ws.Cells.Item(250, 5).Value = 7
$ws.Cells.Item(250, 6).Value = 'Fruta'
$ws.Cells.Item(250, 7).Value = 100108
$ws.Cells.Item(250, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(250, 9).Value = 100108005
$ws.Cells.Item(250, 10).Value = 'Piña'
$ws.Cells.Item(250, 11).Value = 'Caramelo'
$ws.Cells.Item(250, 12).Value = 'Segunda'
$ws.Cells.Item(250, 13).Value = 540
$ws.Cells.Item(250, 14).Value = 19000
$ws.Cells.Item(250, 15).Value = 19000
$ws.Cells.Item(250, 16).Value = 19000
$ws.Cells.Item(250, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(250, 18).Value = 'Ecuador'
$ws.Cells.Item(250, 19).Value = 1357
$ws.Cells.Item(250, 20).Value = 14
